$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.27156410718331
$ws.Range("C2").Value = 9.625922670097504
$ws.Range("D2").Value = 3.575400863928837
$ws.Range("E2").Value = 16.58294073349141
$ws.Range("F2").Value = 20.83286531184558
$ws.Range("N2").Value = 15.75157488290232
$ws.Range("O2").Value = 18.29973775294018
$ws.Range("B3").Value = 12.60019494913002
$ws.Range("C3").Value = 9.134130970340939
$ws.Range("D3").Value = 3.549436474145852
$ws.Range("E3").Value = 15.63564054750957
$ws.Range("F3").Value = 20.70145259169839
$ws.Range("N3").Value = 15.80189921768247
$ws.Range("O3").Value = 18.27781497472204
$ws.Range("B4").Value = 12.17052012960297
$ws.Range("C4").Value = 8.816595894821873
$ws.Range("D4").Value = 3.533264112081607
$ws.Range("E4").Value = 15.02879000762039
$ws.Range("F4").Value = 20.62864604092896
$ws.Range("N4").Value = 15.83465806696178
$ws.Range("O4").Value = 18.27090897552302
$ws.Range("B5").Value = 11.99124259529006
$ws.Range("C5").Value = 8.683362882070044
$ws.Range("D5").Value = 3.526619002687049
$ws.Range("E5").Value = 14.77542321465228
$ws.Range("F5").Value = 20.60098457902129
$ws.Range("N5").Value = 15.84847605398079
$ws.Range("O5").Value = 18.26974068951908
$ws.Range("B6").Value = 11.96122788212026
$ws.Range("C6").Value = 8.661010710042822
$ws.Range("D6").Value = 3.525512372818075
$ws.Range("E6").Value = 14.7329940929246
$ws.Range("F6").Value = 20.59651335638493
$ws.Range("N6").Value = 15.85079884278382
$ws.Range("O6").Value = 18.26964600694011
$ws.Range("B7").Value = 12.16811896194693
$ws.Range("C7").Value = 8.814814476277427
$ws.Range("D7").Value = 3.533174711337349
$ws.Range("E7").Value = 15.02539720001694
$ws.Range("F7").Value = 20.62826482820009
$ws.Range("N7").Value = 15.83484252308212
$ws.Range("O7").Value = 18.27088655959344
$ws.Range("B8").Value = 13.04381155157603
$ws.Range("C8").Value = 9.459628927021098
$ws.Range("D8").Value = 3.566497970554614
$ws.Range("E8").Value = 16.2616876183642
$ws.Range("F8").Value = 20.78593527466846
$ws.Range("N8").Value = 15.76854131182965
$ws.Range("O8").Value = 18.29081754348761
$ws.Range("B9").Value = 14.61513578123397
$ws.Range("C9").Value = 10.597723095996
$ws.Range("D9").Value = 3.629883714108387
$ws.Range("E9").Value = 18.59154783222799
$ws.Range("F9").Value = 21.15632047419537
$ws.Range("N9").Value = 15.6532389882968
$ws.Range("O9").Value = 18.38192759954019
$ws.Range("B10").Value = 15.6726344556512
$ws.Range("C10").Value = 11.35391702605415
$ws.Range("D10").Value = 3.675080027525733
$ws.Range("E10").Value = 20.24807268840485
$ws.Range("F10").Value = 21.46368247961469
$ws.Range("N10").Value = 15.57744171402299
$ws.Range("O10").Value = 18.48048428928114
$ws.Range("B11").Value = 16.13143972783831
$ws.Range("C11").Value = 11.68014848791492
$ws.Range("D11").Value = 3.695305875644059
$ws.Range("E11").Value = 20.95931675908041
$ws.Range("F11").Value = 21.61064990407018
$ws.Range("N11").Value = 15.5448843084209
$ws.Range("O11").Value = 18.53212519522911
$ws.Range("B12").Value = 16.30189895022955
$ws.Range("C12").Value = 11.80110409868711
$ws.Range("D12").Value = 3.702913747811815
$ws.Range("E12").Value = 21.22259862994759
$ws.Range("F12").Value = 21.66728315563908
$ws.Range("N12").Value = 15.53283142729869
$ws.Range("O12").Value = 18.55265092800852
$ws.Range("B13").Value = 16.26533451519838
$ws.Range("C13").Value = 11.77516929915695
$ws.Range("D13").Value = 3.701277584070588
$ws.Range("E13").Value = 21.16616473536934
$ws.Range("F13").Value = 21.65504337810581
$ws.Range("N13").Value = 15.53541497004315
$ws.Range("O13").Value = 18.54818733690996
$ws.Range("B14").Value = 16.14552965839578
$ws.Range("C14").Value = 11.69015143932472
$ws.Range("D14").Value = 3.695932824254285
$ws.Range("E14").Value = 20.98109818846065
$ws.Range("F14").Value = 21.61528980001895
$ws.Range("N14").Value = 15.54388718474472
$ws.Range("O14").Value = 18.53379446519936
$ws.Range("B15").Value = 16.0717164097328
$ws.Range("C15").Value = 11.63773871956856
$ws.Range("D15").Value = 3.692652241612172
$ws.Range("E15").Value = 20.86695260535277
$ws.Range("F15").Value = 21.59106573389948
$ws.Range("N15").Value = 15.54911257279558
$ws.Range("O15").Value = 18.52510451503666
$ws.Range("B16").Value = 15.6421954096912
$ws.Range("C16").Value = 11.33223713789318
$ws.Range("D16").Value = 3.673751224557625
$ws.Range("E16").Value = 20.2007426243093
$ws.Range("F16").Value = 21.45421759340264
$ws.Range("N16").Value = 15.57960805995792
$ws.Range("O16").Value = 18.47724573371858
$ws.Range("B17").Value = 15.37293592415137
$ws.Range("C17").Value = 11.14025079965935
$ws.Range("D17").Value = 3.662068202999523
$ws.Range("E17").Value = 19.78123040252719
$ws.Range("F17").Value = 21.37206294210545
$ws.Range("N17").Value = 15.59880815061481
$ws.Range("O17").Value = 18.44962397603635
$ws.Range("B18").Value = 15.21597382962659
$ws.Range("C18").Value = 11.02815435598203
$ws.Range("D18").Value = 3.655317169047275
$ws.Range("E18").Value = 19.53595291292994
$ws.Range("F18").Value = 21.32548549387339
$ws.Range("N18").Value = 15.61003258807859
$ws.Range("O18").Value = 18.43437799929937
$ws.Range("B19").Value = 15.16247269633603
$ws.Range("C19").Value = 10.98991434969025
$ws.Range("D19").Value = 3.653026111751265
$ws.Range("E19").Value = 19.45222066482444
$ws.Range("F19").Value = 21.30983264517468
$ws.Range("N19").Value = 15.61386410593209
$ws.Range("O19").Value = 18.42932634707763
$ws.Range("B20").Value = 15.40181611141165
$ws.Range("C20").Value = 11.16086122572586
$ws.Range("D20").Value = 3.663315139925944
$ws.Range("E20").Value = 19.82630039271546
$ws.Range("F20").Value = 21.38073884805473
$ws.Range("N20").Value = 15.59674553471742
$ws.Range("O20").Value = 18.45249803935712
$ws.Range("B21").Value = 16.18080883948502
$ws.Range("C21").Value = 11.71519347866379
$ws.Range("D21").Value = 3.697504124719239
$ws.Range("E21").Value = 21.03562066414173
$ws.Range("F21").Value = 21.62694018897354
$ws.Range("N21").Value = 15.54139120701363
$ws.Range("O21").Value = 18.53799573908935
$ws.Range("B22").Value = 16.67077820928148
$ws.Range("C22").Value = 12.06242850951304
$ws.Range("D22").Value = 3.719548642215879
$ws.Range("E22").Value = 21.79072833981449
$ws.Range("F22").Value = 21.79353469548956
$ws.Range("N22").Value = 15.50682184279566
$ws.Range("O22").Value = 18.59952471419356
$ws.Range("B23").Value = 16.41104692461807
$ws.Range("C23").Value = 11.87848722027909
$ws.Range("D23").Value = 3.70781154451398
$ws.Range("E23").Value = 21.39092751753035
$ws.Range("F23").Value = 21.70411604581724
$ws.Range("N23").Value = 15.52512524655821
$ws.Range("O23").Value = 18.56617170026103
$ws.Range("B24").Value = 15.3887661042717
$ws.Range("C24").Value = 11.15154861042544
$ws.Range("D24").Value = 3.6627515063501
$ws.Range("E24").Value = 19.80593700602815
$ws.Range("F24").Value = 21.37681442910545
$ws.Range("N24").Value = 15.59767746377773
$ws.Range("O24").Value = 18.45119669941386
$ws.Range("B25").Value = 14.20659368664019
$ws.Range("C25").Value = 10.30371099249511
$ws.Range("D25").Value = 3.612964612971913
$ws.Range("E25").Value = 17.9437629592119
$ws.Range("F25").Value = 21.04976934090907
$ws.Range("N25").Value = 15.68286184824613
$ws.Range("O25").Value = 18.35171279602193
